$d = $word.ActiveDocument

$d.Content.Find.Execute("144÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "559÷6=", 2) | Out-Null
$d.Content.Find.Execute("150÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "596÷7=", 2) | Out-Null
$d.Content.Find.Execute("426÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "360÷2=", 2) | Out-Null
$d.Content.Find.Execute("776÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "159÷7=", 2) | Out-Null
$d.Content.Find.Execute("316÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "249÷2=", 2) | Out-Null
$d.Content.Find.Execute("924÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "106÷9=", 2) | Out-Null
$d.Content.Find.Execute("899÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "233÷4=", 2) | Out-Null
$d.Content.Find.Execute("482÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "404÷3=", 2) | Out-Null
$d.Content.Find.Execute("620÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "698÷3=", 2) | Out-Null
$d.Content.Find.Execute("880÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "514÷9=", 2) | Out-Null
$d.Content.Find.Execute("592÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "999÷4=", 2) | Out-Null
$d.Content.Find.Execute("756÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "463÷3=", 2) | Out-Null
$d.Content.Find.Execute("847÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "347÷3=", 2) | Out-Null
$d.Content.Find.Execute("368÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "309÷2=", 2) | Out-Null
$d.Content.Find.Execute("769÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "587÷6=", 2) | Out-Null
$d.Content.Find.Execute("224÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "417÷3=", 2) | Out-Null
$d.Content.Find.Execute("118÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "355÷7=", 2) | Out-Null
$d.Content.Find.Execute("851÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "831÷3=", 2) | Out-Null
$d.Content.Find.Execute("114÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "670÷5=", 2) | Out-Null
$d.Content.Find.Execute("637÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "938÷4=", 2) | Out-Null
$d.Content.Find.Execute("922÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "219÷9=", 2) | Out-Null
$d.Content.Find.Execute("203÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "129÷3=", 2) | Out-Null
$d.Content.Find.Execute("425÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "469÷5=", 2) | Out-Null
$d.Content.Find.Execute("182÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "893÷4=", 2) | Out-Null
$d.Content.Find.Execute("605÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "930÷3=", 2) | Out-Null
